# Update "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the a1c53372-a228-4bad-b83b-7164ea0a7679 row (row 4) on each sheet
# to reflect the new report generation run times.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-04 06:49:20"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-04 06:49:16"
$wsZhCn.Range("K4").Value = "2016-09-04 06:49:33"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-04 06:49:20"
$wsDeDe.Range("K4").Value = "2016-09-04 06:49:41"
